# Generate Report for Handback
# Update the handoff/handback timestamps recorded for the
# "c011fe02-d93f-4ee9-93fc-10d3a9508a85" file in both the zh-cn and
# de-de localization report sheets.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E4").Value = "2016-03-21 14:22:41"
$zhcn.Range("E5").Value = "2016-03-21 14:22:41"
$zhcn.Range("H4").Value = "2016-03-21 14:23:05"
$zhcn.Range("H5").Value = "2016-03-21 14:23:05"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E4").Value = "2016-03-21 14:22:45"
$dede.Range("E5").Value = "2016-03-21 14:22:45"
$dede.Range("H4").Value = "2016-03-21 14:23:14"
$dede.Range("H5").Value = "2016-03-21 14:23:14"
